$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 328
$ws.Range("J55").Value = 357.5
$ws.Range("L55").Value = 357.5
$ws.Range("N55").Value = -785.5
$ws.Range("H112").Value = 930.5
$ws.Range("J112").Value = 2222
$ws.Range("L112").Value = 6666
$ws.Range("N112").Value = -8882
$ws.Range("H125").Value = 999.8889
$ws.Range("H137").Value = 1476096.5
$ws.Range("I137").Value = 2085067.5
$ws.Range("J137").Value = 14566.3
$ws.Range("K137").Value = 6255202.5
$ws.Range("L137").Value = 43698.89999999999
$ws.Range("M137").Value = -6252652.5
$ws.Range("N137").Value = -48798.89999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 44596.9
$ws.Range("I45").Value = 141658.33
$ws.Range("K45").Value = 141658.33
$ws.Range("M45").Value = -141281.33
$ws.Range("H61").Value = 3453.2307
$ws.Range("I61").Value = 1974.5
$ws.Range("J61").Value = 3722.0908
$ws.Range("K61").Value = 1974.5
$ws.Range("L61").Value = 3722.0908
$ws.Range("M61").Value = -1762.5
$ws.Range("N61").Value = -4146.0908
$ws.Range("H74").Value = 233694
$ws.Range("I74").Value = 254675.27
$ws.Range("K74").Value = 254675.27
$ws.Range("M74").Value = -253801.27
$ws.Range("H77").Value = 233694
$ws.Range("I77").Value = 254675.27
$ws.Range("K77").Value = 1273376.35
$ws.Range("M77").Value = -1269008.35
$ws.Range("H122").Value = 2455.8845
$ws.Range("I122").Value = 2364.5217
$ws.Range("K122").Value = 7093.5651
$ws.Range("M122").Value = -4643.5651
$ws.Range("H136").Value = 3453.2307
$ws.Range("I136").Value = 1974.5
$ws.Range("J136").Value = 3722.0908
$ws.Range("K136").Value = 5923.5
$ws.Range("L136").Value = 11166.2724
$ws.Range("M136").Value = -3373.5
$ws.Range("N136").Value = -16266.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 1894.5
$ws.Range("I54").Value = 1894.5
$ws.Range("K54").Value = 1894.5
$ws.Range("M54").Value = -1410.5
$ws.Range("H86").Value = 2442.24
$ws.Range("J86").Value = 3499.3333
$ws.Range("L86").Value = 3499.3333
$ws.Range("N86").Value = -5745.3333
$ws.Range("H89").Value = 2442.24
$ws.Range("J89").Value = 3499.3333
$ws.Range("L89").Value = 17496.6665
$ws.Range("N89").Value = -28728.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3935.2827
$ws.Range("I31").Value = 2797.4644
$ws.Range("J31").Value = 5705.222
$ws.Range("K31").Value = 2797.4644
$ws.Range("L31").Value = 5705.222
$ws.Range("M31").Value = -2502.4644
$ws.Range("N31").Value = -6295.222
$ws.Range("H34").Value = 3935.2827
$ws.Range("I34").Value = 2797.4644
$ws.Range("J34").Value = 5705.222
$ws.Range("K34").Value = 2797.4644
$ws.Range("L34").Value = 5705.222
$ws.Range("M34").Value = -2595.4644
$ws.Range("N34").Value = -6109.222
$ws.Range("H58").Value = 2262.85
$ws.Range("I58").Value = 1517.5714
$ws.Range("K58").Value = 1517.5714
$ws.Range("M58").Value = -1314.5714
$ws.Range("H105").Value = 3134.8333
$ws.Range("I105").Value = 3699.75
$ws.Range("J105").Value = 2005
$ws.Range("K105").Value = 3699.75
$ws.Range("L105").Value = 2005
$ws.Range("M105").Value = -1952.75
$ws.Range("N105").Value = -5499
$ws.Range("H122").Value = 2393.1428
$ws.Range("I122").Value = 2688
$ws.Range("K122").Value = 8064
$ws.Range("M122").Value = -5614
$ws.Range("H132").Value = 11909487
$ws.Range("I132").Value = 4390.5835
$ws.Range("K132").Value = 13171.7505
$ws.Range("M132").Value = -10641.7505
$ws.Range("H134").Value = 6968.4
$ws.Range("I134").Value = 7076
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 21228
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -18693
$ws.Range("N134").Value = -23070
$ws.Range("H136").Value = 2262.85
$ws.Range("I136").Value = 1517.5714
$ws.Range("K136").Value = 4552.7142
$ws.Range("M136").Value = -2002.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 10000
$ws.Range("J88").Value = 10000
$ws.Range("L88").Value = 30000
$ws.Range("N88").Value = -30856
$ws.Range("H91").Value = 10000
$ws.Range("J91").Value = 10000
$ws.Range("L91").Value = 30000
$ws.Range("N91").Value = -32964
$ws.Range("H107").Value = 977.25
$ws.Range("I107").Value = 545.8333
$ws.Range("J107").Value = 1408.6666
$ws.Range("K107").Value = 1637.4999
$ws.Range("L107").Value = 4225.9998
$ws.Range("M107").Value = 282.5001
$ws.Range("N107").Value = -8065.9998
$ws.Range("H131").Value = 19233638
$ws.Range("J131").Value = 2610.6667
$ws.Range("L131").Value = 7832.000100000001
$ws.Range("N131").Value = -17912.0001
$ws.Range("H132").Value = 989.4
$ws.Range("J132").Value = 999.5
$ws.Range("L132").Value = 8995.5
$ws.Range("N132").Value = -14055.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6414983
$ws.Range("I122").Value = 10992823
$ws.Range("K122").Value = 32978469
$ws.Range("M122").Value = -32976019
$ws.Range("H132").Value = 3083.111
$ws.Range("I132").Value = 2842.5715
$ws.Range("J132").Value = 3925
$ws.Range("K132").Value = 8527.7145
$ws.Range("L132").Value = 11775
$ws.Range("M132").Value = -5997.7145
$ws.Range("N132").Value = -16835

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12142.714
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5272
$ws.Range("H136").Value = 5133.5884
$ws.Range("I136").Value = 3797
$ws.Range("J136").Value = 6321.6665
$ws.Range("K136").Value = 11391
$ws.Range("L136").Value = 18964.9995
$ws.Range("M136").Value = -8841
$ws.Range("N136").Value = -24064.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4228.2
$ws.Range("I81").Value = 3921.077
$ws.Range("J81").Value = 4798.5713
$ws.Range("K81").Value = 7842.154
$ws.Range("L81").Value = 9597.142599999999
$ws.Range("M81").Value = -6781.154
$ws.Range("N81").Value = -11719.1426
$ws.Range("H84").Value = 4228.2
$ws.Range("I84").Value = 3921.077
$ws.Range("J84").Value = 4798.5713
$ws.Range("K84").Value = 39210.77
$ws.Range("L84").Value = 47985.713
$ws.Range("M84").Value = -33906.77
$ws.Range("N84").Value = -58593.713
$ws.Range("H122").Value = 19235544
$ws.Range("I122").Value = 5388.273
$ws.Range("J122").Value = 125001400
$ws.Range("K122").Value = 16164.819
$ws.Range("L122").Value = 375004200
$ws.Range("M122").Value = -13714.819
$ws.Range("N122").Value = -375009100
$ws.Range("H126").Value = 2863.7693
$ws.Range("I126").Value = 2892.7144
$ws.Range("K126").Value = 8678.143199999999
$ws.Range("M126").Value = -6208.143199999999
$ws.Range("H132").Value = 1432.7073
$ws.Range("I132").Value = 1395.9231
$ws.Range("J132").Value = 2150
$ws.Range("K132").Value = 4187.7693
$ws.Range("L132").Value = 6450
$ws.Range("M132").Value = -1657.7693
$ws.Range("N132").Value = -11510
$ws.Range("H136").Value = 583204.75
$ws.Range("I136").Value = 3518.875
$ws.Range("J136").Value = 2129033.8
$ws.Range("K136").Value = 10556.625
$ws.Range("L136").Value = 6387101.399999999
$ws.Range("M136").Value = -8006.625
$ws.Range("N136").Value = -6392201.399999999
